# Lab 4 SAWA tables - update Table 4.1 measurements
# - clear the (now unused) "Guessed_Distance" column D for rows 3-13
# - correct several Team-1 (F) and Instructor (G) measurements
# Table 4.2's formulas pull from Table 4.1!H column, so they (and the
# chart caches) recompute automatically once the workbook recalculates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 4.1")

# Clear column D (Guessed_Distance) for all data rows
$ws.Range("D3:D13").ClearContents()

# Row 3 - Base of Closest Magnolia Tree (E3 stays 39, F3 stays 41)
$ws.Range("G3").Value = 15

# Row 4 - Lamp Post across sidewalk
$ws.Range("F4").Value = 25
$ws.Range("G4").Value = 12

# Row 5 - Lamp Post on Right... (E5 stays 16, F5 stays 16)
$ws.Range("G5").Value = 10

# Row 6 - Base of Brick Wall
$ws.Range("F6").Value = 84
$ws.Range("G6").Value = 20

# Row 7 - Red fire hydrant to right
$ws.Range("F7").Value = 143
$ws.Range("G7").Value = 127

# Row 8 - Shortest distance to letters on Sidewalk
$ws.Range("F8").Value = 277
$ws.Range("G8").Value = 261

# Row 9 - Nearest point to manhole cover to left
$ws.Range("F9").Value = 9
$ws.Range("G9").Value = 7

# Row 10 - Closest corner to concrete square to right
$ws.Range("F10").Value = 42
$ws.Range("G10").Value = 42

# Row 11 - Facing statue first tree past oval walkway
$ws.Range("F11").Value = 107
$ws.Range("G11").Value = 107

# Row 12 - Leftmost urn
$ws.Range("F12").Value = 196
$ws.Range("G12").Value = 196

# Update the sheet view/selection to match the reviewed range
$ws.Activate()
$ws.Range("H3:H12").Select()
$excel.ActiveWindow.ScrollColumn = 5
